# Rewrites A1 with the reformatted (pretty-printed, double-quoted JSON)
# questions payload and drops the stray A2 cell / leftover bold+bordered
# style that the old single-line layout used.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newText = "questions = [`n    {`n        ""title"": ""You want to implement a frontend library to handle GraphQL subscriptions. At first, you try to implement it using HTTP requests, but you find that this requires you to continuously launch new requests to check for updates.What kind of connection should you use instead?"",`n        ""ques_type"": 2,`n        ""options"": [`n            ""WebSockets"",`n            ""WebAssembly"",`n            ""AJAX requests"",`n            ""REST requests""`n        ],`n        ""score"": ""WebSockets""`n    },`n    {`n        ""title"": ""You are creating an API for a pizzeria that sells pizza, soda, and bread. You want to use a single type in GraphQL that represents all three of these types of food using an enum type.What syntax should you use to accomplish this?"",`n        ""ques_type"": 2,`n        ""options"": [`n            ""enum Order {\n Pizza, Soda, Bread\n}\n"",`n            ""ENUM ORDER {\n  Pizza, Soda, Bread\n}\n"",`n            ""enum Order {\n PIZZA\n SODA\n BREAD\n}\n"",`n            ""Enumerate order { pizza soda bread }\n""`n        ],`n        ""score"": ""enum Order {\n PIZZA\n SODA\n BREAD\n}""`n    },`n    {`n        ""title"": ""You are developing a GraphQL API with Apollo GraphQL that fetches data from an external REST API. Occasionally, the REST API returns a 500 status code due to internal server errors. How should your resolver gracefully handle this error?"",`n        ""ques_type"": 2,`n        ""options"": [`n            ""Use a try-catch block in the resolver function to catch the error and return a custom error message to the client."",`n            ""Retry the request automatically in the resolver function twice before returning an error to the client."",`n            ""Return a generic 500 internal server error message and empty JSON object to the client."",`n            ""Use Apollo Server's createError function to generate a GraphQL error object and return it to the client.""`n        ],`n        ""score"": ""Use Apollo Server's createError function to generate a GraphQL error object and return it to the client.""`n    },`n    {`n        ""title"": ""Attackers are targeting your GraphQL API as part of a denial of service (DoS) attack. They are using a botnet to apply various different DoS techniques simultaneously.How should you mitigate this threat?"",`n        ""ques_type"": 15,`n        ""options"": [`n            ""Block all UDP traffic to your API."",`n            ""Disable requests that come from multiple IP addresses."",`n            ""Use query cost analysis."",`n            ""Implement query depth limiting."",`n            ""Ban IPs that make more than one request per minute."",`n            ""Set a timeout in the HTTP server configuration.""`n        ],`n        ""score"": [`n            ""Use query cost analysis."",`n            ""Implement query depth limiting."",`n            ""Set a timeout in the HTTP server configuration.""`n        ]`n    }`n]"

# A2 held the old long string underneath A1's numeric placeholder; delete
# the whole row so the sheet collapses back down to just A1.
$ws.Rows(2).Delete()

# A1 used to be styled bold/centered/bordered (s="1"); clear that back to
# the workbook's default formatting before writing the new text.
$ws.Range("A1").ClearFormats()
$ws.Range("A1").Value = $newText

# Ensure no stray auto-calculated row height sticks around now that A1
# contains a many-line string.
$ws.Rows(1).AutoFit()
